# Apply updated loading-percent results ("case with 380 kV done") to Sheet1.
# The table data spans rows 2-25 (data index 0-23) and columns B,C,D,E,G,H,I,L,M
# (columns F, J, K, N, O are unchanged zero placeholders).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry: row, B, C, D, E, G, H, I, L, M
$data = @(
    @(2, 17.2155782130025, 9.265561133906502, 6.012369528219663, 10.72001598602268, 47.36080490319004, 18.50961994386905, 28.44064089923697, 10.23932118325357, 15.93095022874076),
    @(3, 16.79327670044246, 8.738709156970289, 5.899378627280349, 10.72313299291434, 47.01949507957875, 18.51745687564921, 28.47476958243126, 10.25008044792773, 15.85688236496242),
    @(4, 16.53297367771291, 8.396896010779011, 5.830835486753257, 10.72527031298915, 46.82533927364572, 18.52662893710382, 28.50269333588268, 10.2582347021799, 15.81453004577969),
    @(5, 16.42681468023834, 8.253037080249044, 5.803154326322792, 10.726197546125, 46.75015664519585, 18.53145990719149, 28.51581703704111, 10.2619469057873, 15.79807029809476),
    @(6, 16.4091871133401, 8.228874958317315, 5.798574233288628, 10.72635491230507, 46.7379120980879, 18.53232801581084, 28.51810134328502, 10.26258682493655, 15.79538580259596),
    @(7, 16.53154208877649, 8.394974310344191, 5.830461099493013, 10.72528259008966, 46.82430931580963, 18.5266896672267, 28.50286327402099, 10.25828319011773, 15.81430481144125),
    @(8, 17.07027795394641, 9.087726688622116, 5.973261852254183, 10.72104438818343, 47.23996304933378, 18.51141545156352, 28.45095766504734, 10.24270971383151, 15.90477112288103),
    @(9, 18.11161663913006, 10.29954771949251, 6.258114660993567, 10.71450377759007, 48.17418200662077, 18.51617595662076, 28.40478241164735, 10.22445071041644, 16.10635088961836),
    @(10, 18.85851447391293, 11.09932437679344, 6.467977660882647, 10.71077464972504, 48.9284320891031, 18.54096359784569, 28.40515283418767, 10.2185167648848, 16.26829109680886),
    @(11, 19.19264208369724, 11.44338530885853, 6.563098626168819, 10.70931130065085, 49.2851510365814, 18.55687827068272, 28.41283637215262, 10.21743900230565, 16.34475799606569),
    @(12, 19.31822763174454, 11.57082598523958, 6.59903087730375, 10.70879063621854, 49.42208327573625, 18.56357189938204, 28.41683070770708, 10.21726367442215, 16.37409780968713),
    @(13, 19.29122436838442, 11.543506067755, 6.591296725122977, 10.7089012824287, 49.39251190134267, 18.56210064367216, 28.41592215575377, 10.21729108780144, 16.36776221223549),
    @(14, 19.20299365780049, 11.45392702480097, 6.566056769469282, 10.70926779461003, 49.29638007050492, 18.55741559207686, 28.41314322826553, 10.21741991490256, 16.34716422675794),
    @(15, 19.14882351934785, 11.39868629731411, 6.550584008157635, 10.70949665214433, 49.2377343158054, 18.55463272790584, 28.41158242756025, 10.21752912930157, 16.33459672579434),
    @(16, 18.83655358514296, 11.07644093274188, 6.461751003467375, 10.71087496925235, 48.90538450133624, 18.54001689733745, 28.40480228837834, 10.21861979898322, 16.2633485194837),
    @(17, 18.64344559669889, 10.87368695172762, 6.40713831918453, 10.71178018365154, 48.70490963495114, 18.53223886993806, 28.40257086054532, 10.2197039749854, 16.22034368251675),
    @(18, 18.53185127132874, 10.75520828129751, 6.375696080624008, 10.71232277687881, 48.59088919941023, 18.52820180331174, 28.40199456253956, 10.22048022752345, 16.19587312572087),
    @(19, 18.49398153966661, 10.71477423758458, 6.365046237507708, 10.71251025868789, 48.552508111038, 18.52690990460724, 28.40192075431686, 10.22076928379719, 16.18763384754002),
    @(20, 18.66405745111936, 10.89546301133519, 6.412955367825409, 10.71168155181655, 48.72611804010476, 18.53302166048807, 28.40273517929724, 10.21957276511566, 16.22489436306383),
    @(21, 19.22893569116273, 11.4803158783466, 6.573473039134117, 10.7091592329098, 49.32456696427784, 18.55877360345537, 28.41392999861476, 10.21737576073222, 16.35320410278623),
    @(22, 19.59257629062784, 11.84595185131055, 6.677851367605042, 10.70770584895654, 49.726422055769, 18.57949185763827, 28.42757047274429, 10.21729657170809, 16.43928774164049),
    @(23, 19.3990417684854, 11.65232498588933, 6.622203228695057, 10.70846370769591, 49.51099839513898, 18.56807854180503, 28.41971052092451, 10.21721485078747, 16.39314609020026),
    @(24, 18.65474060949698, 10.88562401205327, 6.410325616410592, 10.7117260742117, 48.71652586447856, 18.53266640703193, 28.40265869013431, 10.2196316086669, 16.22283620932342),
    @(25, 17.83248548237097, 9.98756251219223, 6.180785506668854, 10.71608397229829, 47.90915989641322, 18.51115624859378, 28.41127677673567, 10.22807580400953, 16.0493192908594)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $ws.Range("D$row").Value = $entry[3]
    $ws.Range("E$row").Value = $entry[4]
    $ws.Range("G$row").Value = $entry[5]
    $ws.Range("H$row").Value = $entry[6]
    $ws.Range("I$row").Value = $entry[7]
    $ws.Range("L$row").Value = $entry[8]
    $ws.Range("M$row").Value = $entry[9]
}
